$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Capture the formatting of D72 (style s="22") BEFORE it gets
#     re-styled below, so the new D76 cell can reuse that same xf. ---
$ws.Range("D72").Copy()
$ws.Range("D76").PasteSpecial(-4122)

# --- Re-style D72:D75 from s="22" to s="18" (copy format from D70,
#     which already carries style s="18"). ---
$ws.Range("D70").Copy()
$ws.Range("D72:D75").PasteSpecial(-4122)

# --- New row 76 ---
# A76: TCID, same look as A75 (s="1")
$ws.Range("A75").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$ws.Range("A76").Value = "Authoring75"

# B76: JIRA id list, same look as B72 (s="8")
$ws.Range("B72").Copy()
$ws.Range("B76").PasteSpecial(-4122)
$ws.Range("B76").Value = "OPQA-1195|OPQA-1313|OPQA-1312|OPQA-1090|OPQA-1201"

# C76: Description, same look as C75 (s="1")
$ws.Range("C75").Copy()
$ws.Range("C76").PasteSpecial(-4122)
$ws.Range("C76").Value = "Verify saving post as draft, accessing it for edit from profile,delete post from prfile"

# D76 value (format already copied above, from the pre-restyle D72)
$ws.Range("D76").Value = "Y"

# E76: empty cell, same look as E75 (s="1")
$ws.Range("E75").Copy()
$ws.Range("E76").PasteSpecial(-4122)

# --- New row 77 ---
# A77: TCID, same look as A75 (s="1")
$ws.Range("A75").Copy()
$ws.Range("A77").PasteSpecial(-4122)
$ws.Range("A77").Value = "Authoring76"

# B77: JIRA id list, plain look (s="1"), copy from B75's column-neutral peer
$ws.Range("A75").Copy()
$ws.Range("B77").PasteSpecial(-4122)
$ws.Range("B77").Value = "OPQA-1196|OPQA-1200|OPQA-1199"

# C77: Description, same look as C75 (s="1")
$ws.Range("C75").Copy()
$ws.Range("C77").PasteSpecial(-4122)
$ws.Range("C77").Value = "Verify draft title,access and edit draft post from post modal, delete post from post modal"

# D77: plain look (s="1")
$ws.Range("A75").Copy()
$ws.Range("D77").PasteSpecial(-4122)
$ws.Range("D77").Value = "Y"

# E77: empty cell, same look as E75 (s="1")
$ws.Range("E75").Copy()
$ws.Range("E77").PasteSpecial(-4122)

# --- Selection / scroll position update ---
$ws.Activate()
$ws.Range("D2:D75").Select()
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
